$wb = $excel.ActiveWorkbook

# Helper: write a value as literal text (not auto-coerced to a number),
# without leaving a residual number-format style override on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. Insert a new sheet "2022-Q1" right before "总计" (the last sheet)
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($zongji)
$newSheet.Name = "2022-Q1"

# Copy the header-row formatting (bold / border / centered) from an existing
# per-fund sheet (2021-Q4) so the new sheet matches the established style.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats
$template.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122) # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (A = running index, B = code, C = name, D..G = text numbers, H = rank)
$rows = @(
    @(0, "008099", "广发价值领先混合",         "61.82", "83.88", "4.28", "2.6459", 7),
    @(1, "270001", "广发聚富混合",             "19.90", "73.54", "4.81", "0.9572", 7),
    @(2, "001763", "广发多策略灵活配置混合",     "20.27", "69.36", "3.64", "0.7378", 9),
    @(3, "002270", "东吴安盈量化灵活配置混合",   "5.24",  "46.02", "2.34", "0.1226", 9),
    @(4, "290012", "泰信行业精选灵活配置混合A",  "0.76",  "92.62", "5.48", "0.0416", 7),
    @(5, "159855", "银华中证影视主题ETF",       "0.96",  "97.27", "3.86", "0.0371", 9),
    @(6, "516620", "国泰中证影视主题ETF",       "0.33",  "96.08", "3.81", "0.0126", 9),
    @(7, "002583", "泰信行业精选灵活配置混合C",  "0.00",  "92.62", "5.48", $null,    7)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newSheet.Range("A$r").Value = $row[0]
    Set-TextValue $newSheet.Range("B$r") $row[1]
    Set-TextValue $newSheet.Range("C$r") $row[2]
    Set-TextValue $newSheet.Range("D$r") $row[3]
    Set-TextValue $newSheet.Range("E$r") $row[4]
    Set-TextValue $newSheet.Range("F$r") $row[5]
    if ($null -eq $row[6]) {
        $newSheet.Range("G$r").Value = 0
    } else {
        Set-TextValue $newSheet.Range("G$r") $row[6]
    }
    $newSheet.Range("H$r").Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing quarters down by one row (values rewritten directly to avoid
#    any Insert()-related formatting/precision artifacts).
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @("2022-Q1", 8,  4.55),
    @("2021-Q4", 15, 4.86),
    @("2021-Q3", 8,  6.31),
    @("2021-Q2", 5,  1.82),
    @("2021-Q1", 2,  0.03),
    @("2020-Q4", 6,  3.48)
)

# Make sure column A has the established bold/border/centered style as far
# down as the new last row (row 7).
$summarySheet.Range("A2").Copy()
$summarySheet.Range("A2:A7").PasteSpecial(-4122) # xlPasteFormats

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summarySheet.Range("A$r").Value = $i
    $summarySheet.Range("B$r").Value = $row[0]
    $summarySheet.Range("C$r").Value = $row[1]
    $summarySheet.Range("D$r").Value = $row[2]
}
